# Apply the "fix: complete task and event data" edit to the imagerating log workbook.
#
# Summary of the change:
#  - Column B is widened slightly (Type column).
#  - Rows 2-11 (previously blank placeholder rows with 0/blank/0/0) are filled in
#    with the first 10 trials worth of Trial/Type data; row 12's Trial number is
#    bumped to 11 to continue the sequence, and rows 2-11's Response/ReactionTime
#    columns are cleared out (no longer 0, just blank) because those trials have
#    not been responded to yet.
#  - Rows 79-81, which previously held the "overflow" trial data that actually
#    belongs near the top of the sheet, are reset back to the blank placeholder
#    state (0 / blank / 0 / 0) now that the data has been moved up to rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (Type) a bit.
$ws.Columns.Item(2).ColumnWidth = 7.6666666666666667

# Fill in the Trial/Type data for rows 2-11, and blank out Response/ReactionTime.
$newData = @(
    @(2, 1, "R7-9.jpg"),
    @(3, 2, "R6-2.jpg"),
    @(4, 3, "T4-6.jpg"),
    @(5, 4, "T3-3.jpg"),
    @(6, 5, "R1-8.jpg"),
    @(7, 6, "T4-3.jpg"),
    @(8, 7, "T5-8.jpg"),
    @(9, 8, "R8-1.jpg"),
    @(10, 9, "T2-4.jpg"),
    @(11, 10, "T2-10.jpg")
)

foreach ($entry in $newData) {
    $row = $entry[0]
    $trial = $entry[1]
    $type = $entry[2]

    $ws.Cells.Item($row, 1).Value = $trial
    $ws.Cells.Item($row, 2).Value = $type
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
}

# Row 12 continues the Trial sequence; Type/Response/ReactionTime stay blank.
$ws.Cells.Item(12, 1).Value = 11

# Rows 79-81 go back to being unused placeholder rows.
for ($row = 79; $row -le 81; $row++) {
    $ws.Cells.Item($row, 1).Value = 0
    $ws.Cells.Item($row, 2).Value = ""
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
